# Add summary info to the "OODP" sheet (screening funnel numbers + notes
# about articles that were not found in the source dataset).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OODP")
$ws.Activate()

# Row 3 - "After screening": not applicable counts + warning note (red text)
$ws.Range("B3").Value = "N.A."
$ws.Range("C3").Value = "N.A."
$ws.Range("D3").Value = "N.A."
$ws.Range("E3").Value = "144 supposed according to article, but not in source excel"
$ws.Range("E3").Font.Color = 255

# Row 2 - "Total": Source / Not in dataset / Final counts + note
$ws.Range("B2").Value = 685
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 685
$ws.Range("E2").Value = "query results"

# Row 4 - "After snowballing": not applicable counts
$ws.Range("B4").Value = "N.A."
$ws.Range("C4").Value = "N.A."
$ws.Range("D4").Value = "N.A."

# Row 5 - "Final selection": Source / Not in dataset / Final counts + note
$ws.Range("B5").Value = 34
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = "Included"

# Row 7 - explanatory note (red text) ahead of the "Not in dataset articles" table
$ws.Range("A7").Value = "The following articles were not found, but some metadata were already in data source such as the abstract, so they were kept."
$ws.Range("A7").Font.Color = 255

# Move the selection/view to D4 (matches the author's final cursor position)
$ws.Range("D4").Select()
